# "Generate Report for Handoff" - update the localization-status report:
#   - Status flips from "In Translation" to "Ready for handoff" (Overview!E2:F2,
#     zh-cn!C2, de-de!C2 all shared the same underlying string).
#   - The two timestamp cells that recorded that status change move forward by
#     ~35-36s (Overview!G2 / de-de!H2 share one timestamp string, zh-cn!H2 has
#     its own).
#   - The now-wider "Status"/"zh-cn"/"de-de" columns are widened to fit the
#     longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps bumped to the handoff-generation time ---
$wsOverview.Range("G2").Value = "2016-08-29 00:38:49"
$wsDeDe.Range("H2").Value = "2016-08-29 00:38:49"
$wsZhCn.Range("H2").Value = "2016-08-29 00:38:45"

# --- Widen the Status-ish columns to fit "Ready for handoff" ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
